# Weekly update: insert a new daily price record at row 19 (Poroto verde,
# Terminal Hortofrutícola Agro Chillán), pushing the previously existing
# rows 19-48 down to 20-49.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 19, shifting rows 19:48 -> 20:49
$ws.Rows("19:19").Insert()

# Populate the newly inserted row 19 with the new record
$ws.Range("A19").Value = 7
$ws.Range("B19").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C19").Value = "Ñuble"
$ws.Range("D19").Value = 44546
$ws.Range("E19").Value = 16
$ws.Range("F19").Value = 100112031
$ws.Range("G19").Value = "Poroto verde"
$ws.Range("H19").Value = "Sin especificar"
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 120
$ws.Range("K19").Value = 12500
$ws.Range("L19").Value = 13000
$ws.Range("M19").Value = 12750
$ws.Range("N19").Value = "$/saco 25 kilos"
$ws.Range("O19").Value = "Región del Maule"
$ws.Range("P19").Value = 510
$ws.Range("Q19").Value = 25
$ws.Range("R19").Value = "Hortaliza"
